# Insert a new data row right above the current row 130 (Fecha = 2021-11-22 / serial 44522),
# pushing the existing rows 130-152 down to 131-153.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(130).EntireRow.Insert()

$ws.Cells.Item(130, 1).Value  = 7
$ws.Cells.Item(130, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(130, 3).Value  = "Ñuble"
$ws.Cells.Item(130, 4).Value  = 44522
$ws.Cells.Item(130, 5).Value  = 16
$ws.Cells.Item(130, 6).Value  = 100112006
$ws.Cells.Item(130, 7).Value  = "Repollo"
$ws.Cells.Item(130, 8).Value  = "Crespo record"
$ws.Cells.Item(130, 9).Value  = "Primera"
$ws.Cells.Item(130, 10).Value = 300
$ws.Cells.Item(130, 11).Value = 600
$ws.Cells.Item(130, 12).Value = 700
$ws.Cells.Item(130, 13).Value = 650
$ws.Cells.Item(130, 14).Value = "$/unidad"
$ws.Cells.Item(130, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(130, 16).Value = 650
$ws.Cells.Item(130, 17).Value = 1
$ws.Cells.Item(130, 18).Value = "Hortaliza"
